# Generate Report for Handback
#
# The "ab33e234-c250-4f33-9a6d-87ebf8d2dd54.md" file has just been handed
# back (in sync with en-US). Update its status on all three sheets and
# stamp the handback datetime on the two locale sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: both locale status columns for the ab33e234 row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn sheet: status + handback datetime for the ab33e234 row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $status
$zhcn.Range("G3").Value = "2016-01-28 04:07:14"

# --- de-de sheet: status + handback datetime for the ab33e234 row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $status
$dede.Range("G3").Value = "2016-01-28 04:07:31"
